# Update vertex-map coordinate labels for the newly-added 2nd graph level.
# Most "V_x_0" (left column) labels shift x from -586 to -583, and the
# 2nd-level row ("V_1_*", y=-143) labels shift y from -143 to -137.
# One label ("V_0_12" -> "( 350, -252)") had its coordinate text split
# across two runs; it gets consolidated into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, [int]$id) {
  for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq $id) {
      return $shp
    }
  }
  throw "Shape with id $id not found"
}

# Replace a (unique) substring of a shape's text with new text, preserving
# the run formatting of whichever run(s) originally held it.
function Set-CoordText($shape, [string]$oldText, [string]$newText) {
  $tr = $shape.TextFrame.TextRange
  $full = $tr.Text
  $idx = $full.IndexOf($oldText)
  if ($idx -lt 0) {
    throw ("Substring [{0}] not found in shape {1} text [{2}]" -f $oldText, $shape.Id, $full)
  }
  $rng = $tr.Characters($idx + 1, $oldText.Length)
  $rng.Text = $newText
}

# id -> (old coordinate text, new coordinate text)
$changes = @(
  @{ Id = 8;   Old = "( -586, -258)"; New = "( -583, -258)" },
  @{ Id = 9;   Old = "( -586, -143)"; New = "( -583, -143)" },
  @{ Id = 10;  Old = "( -586, -26)";  New = "( -583, -26)" },
  @{ Id = 11;  Old = "( -586, 86)";   New = "( -583, 86)" },
  @{ Id = 27;  Old = "( -586, 198)";  New = "( -583, 198)" },
  @{ Id = 87;  Old = "( -586, 273)";  New = "( -583, 273)" },

  @{ Id = 104; Old = "( -486, -143)"; New = "( -486, -137)" },
  @{ Id = 106; Old = "( -586, -252)"; New = "( -583, -252)" },
  @{ Id = 107; Old = "( -586, -143)"; New = "( -583, -137)" },
  @{ Id = 108; Old = "( -586, -26)";  New = "( -583, -26)" },
  @{ Id = 109; Old = "( -586, 86)";   New = "( -583, 86)" },
  @{ Id = 110; Old = "( -426, -143)"; New = "( -426, -137)" },
  @{ Id = 111; Old = "( -306, -143)"; New = "( -306, -137)" },
  @{ Id = 112; Old = "( -380, -143)"; New = "( -380, -137)" },
  @{ Id = 113; Old = "( -253, -143)"; New = "( -253, -137)" },
  @{ Id = 114; Old = "( -130, -143)"; New = "( -130, -137)" },
  @{ Id = 115; Old = "( -186, -143)"; New = "( -186, -137)" },
  @{ Id = 116; Old = "( 50, -143)";   New = "( 50, -137)" },
  @{ Id = 117; Old = "( 173, -143)";  New = "( 173, -137)" },
  @{ Id = 118; Old = "( 106, -143)";  New = "( 106, -137)" },
  @{ Id = 119; Old = "( 233, -143)";  New = "( 233, -137)" },
  @{ Id = 120; Old = "( 350, -143)";  New = "( 350, -137)" },
  @{ Id = 121; Old = "( 293, -143)";  New = "( 293, -137)" },
  @{ Id = 122; Old = "( 406, -143)";  New = "( 406, -137)" },
  @{ Id = 123; Old = "( 466, -143)";  New = "( 466, -137)" },
  @{ Id = 124; Old = "( 526, -143)";  New = "( 526, -137)" },
  @{ Id = 125; Old = "( -586, 198)";  New = "( -583, 198)" }
)

foreach ($chg in $changes) {
  $shape = Get-ShapeById $s $chg.Id
  Set-CoordText $shape $chg.Old $chg.New
}

# Shape id 135 ("TextBox 134", label V_0_12) stores its coordinate text
# split across two runs: "( 350, " + "-252)". Consolidate into one run
# reading "( 350, -252)" while keeping the first run's formatting.
$shape135 = Get-ShapeById $s 135
$tr135 = $shape135.TextFrame.TextRange
$full135 = $tr135.Text
$tailOld = "-252)"
$tailIdx = $full135.IndexOf($tailOld)
if ($tailIdx -lt 0) {
  throw ("Tail substring not found in shape 135 text [{0}]" -f $full135)
}
$tailRng = $tr135.Characters($tailIdx + 1, $tailOld.Length)
$tailRng.Delete()

$headOld = "( 350, "
$full135b = $tr135.Text
$headIdx = $full135b.IndexOf($headOld)
if ($headIdx -lt 0) {
  throw ("Head substring not found in shape 135 text [{0}]" -f $full135b)
}
$headRng = $tr135.Characters($headIdx + 1, $headOld.Length)
$headRng.Text = "( 350, -252)"
